# Applies the trade-close + 3 new OPEN trades update described in the commit:
#   "Trade #104 closed at 2026-02-18 00:37:41 - unknown UNKNOWN +0.000%"
#
# Affected sheets: Summary, Strategy Status, All Trades, momentum,
#                   HighProbConvergence, MarketMaking

# Helper: write a text value to a cell while preventing Excel's automatic
# "looks like a date" (or other) conversion from changing the stored type,
# and without leaving an explicit number-format style applied to the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - roll-up numbers after the trade close + new trades
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(3, 2).Value = 1499.06   # Current Capital
$wsSummary.Cells.Item(4, 2).Value = 0.17      # Total P&L $
$wsSummary.Cells.Item(6, 2).Value = 132       # Total Trades
$wsSummary.Cells.Item(8, 2).Value = 47        # Losing Trades
$wsSummary.Cells.Item(9, 2).Value = 46.97     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Cells.Item(6, 3).Value = 99.14      # Capital
$wsStatus.Cells.Item(6, 4).Value = 52         # Trades
$wsStatus.Cells.Item(6, 5).Value = -0.67      # P&L $
$wsStatus.Cells.Item(6, 6).Value = -0.86      # P&L %
$wsStatus.Cells.Item(6, 7).Value = 46.15      # Win Rate %

# ---------------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# 3a) Close trade #132 (row 133, MarketMaking DOWN)
$wsAll.Cells.Item(133, 7).Value = 0.03                # Exit Price
$wsAll.Cells.Item(133, 8).Value = "CLOSED"            # Status
$wsAll.Cells.Item(133, 9).Value = -25                 # P&L %
$wsAll.Cells.Item(133, 10).Value = -0.01              # P&L $
$wsAll.Cells.Item(133, 11).Value = 99.14              # Capital After
$wsAll.Cells.Item(133, 12).Value = "early_exit"       # Exit Reason
$wsAll.Cells.Item(133, 13).Value = 0.15               # Duration (min)

# 3b) New row 162 - trade #161 (momentum DOWN, OPEN)
$wsAll.Cells.Item(162, 1).Value = 161
Set-TextValue $wsAll.Cells.Item(162, 2) "2026-02-18"
Set-TextValue $wsAll.Cells.Item(162, 3) "00:37:34"
$wsAll.Cells.Item(162, 4).Value = "momentum"
$wsAll.Cells.Item(162, 5).Value = "DOWN"
$wsAll.Cells.Item(162, 6).Value = 0.04
$wsAll.Cells.Item(162, 8).Value = "OPEN"
$wsAll.Cells.Item(162, 9).Value = 0
$wsAll.Cells.Item(162, 10).Value = 0
$wsAll.Cells.Item(162, 11).Value = 99.23374292899115
$wsAll.Cells.Item(162, 13).Value = 0
$wsAll.Cells.Item(162, 14).Value = 0
$wsAll.Cells.Item(162, 15).Value = 0
$wsAll.Cells.Item(162, 16).Value = 0.9
$wsAll.Cells.Item(162, 17).Value = "Downward momentum: -1.942% over 10 samples"

# 3c) New row 163 - trade #162 (HighProbConvergence UP, OPEN)
$wsAll.Cells.Item(163, 1).Value = 162
Set-TextValue $wsAll.Cells.Item(163, 2) "2026-02-18"
Set-TextValue $wsAll.Cells.Item(163, 3) "00:37:34"
$wsAll.Cells.Item(163, 4).Value = "HighProbConvergence"
$wsAll.Cells.Item(163, 5).Value = "UP"
$wsAll.Cells.Item(163, 6).Value = 0.96
$wsAll.Cells.Item(163, 8).Value = "OPEN"
$wsAll.Cells.Item(163, 9).Value = 0
$wsAll.Cells.Item(163, 10).Value = 0
$wsAll.Cells.Item(163, 11).Value = 100.4130057263667
$wsAll.Cells.Item(163, 13).Value = 0
$wsAll.Cells.Item(163, 14).Value = 0
$wsAll.Cells.Item(163, 15).Value = 0
$wsAll.Cells.Item(163, 16).Value = 0.95
$wsAll.Cells.Item(163, 17).Value = "Mean reversion UP: price 1.75% below mean (z=-3.00)"

# 3d) New row 164 - trade #163 (MarketMaking UP, OPEN)
$wsAll.Cells.Item(164, 1).Value = 163
Set-TextValue $wsAll.Cells.Item(164, 2) "2026-02-18"
Set-TextValue $wsAll.Cells.Item(164, 3) "00:37:35"
$wsAll.Cells.Item(164, 4).Value = "MarketMaking"
$wsAll.Cells.Item(164, 5).Value = "UP"
$wsAll.Cells.Item(164, 6).Value = 0.9399999999999999
$wsAll.Cells.Item(164, 8).Value = "OPEN"
$wsAll.Cells.Item(164, 9).Value = 0
$wsAll.Cells.Item(164, 10).Value = 0
$wsAll.Cells.Item(164, 11).Value = 99.14858346467945
$wsAll.Cells.Item(164, 13).Value = 0
$wsAll.Cells.Item(164, 14).Value = 0
$wsAll.Cells.Item(164, 15).Value = 0
$wsAll.Cells.Item(164, 16).Value = 0.6
$wsAll.Cells.Item(164, 17).Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# 4) momentum sheet - append new row 42 (trade #161, DOWN, OPEN)
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Cells.Item(42, 1).Value = 161
Set-TextValue $wsMomentum.Cells.Item(42, 2) "2026-02-18"
Set-TextValue $wsMomentum.Cells.Item(42, 3) "00:37:34"
$wsMomentum.Cells.Item(42, 4).Value = "momentum"
$wsMomentum.Cells.Item(42, 5).Value = "DOWN"
$wsMomentum.Cells.Item(42, 6).Value = 0.04
$wsMomentum.Cells.Item(42, 8).Value = "OPEN"
$wsMomentum.Cells.Item(42, 9).Value = 0
$wsMomentum.Cells.Item(42, 10).Value = 0
$wsMomentum.Cells.Item(42, 11).Value = 99.23374292899115
$wsMomentum.Cells.Item(42, 12).Value = 0
$wsMomentum.Cells.Item(42, 13).Value = 0
$wsMomentum.Cells.Item(42, 14).Value = 0.9
$wsMomentum.Cells.Item(42, 15).Value = "Downward momentum: -1.942% over 10 samples"
$wsMomentum.Cells.Item(42, 17).Value = 0

# ---------------------------------------------------------------------------
# 5) HighProbConvergence sheet - append new row 24 (trade #162, UP, OPEN)
# ---------------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
$wsHPC.Cells.Item(24, 1).Value = 162
Set-TextValue $wsHPC.Cells.Item(24, 2) "2026-02-18"
Set-TextValue $wsHPC.Cells.Item(24, 3) "00:37:34"
$wsHPC.Cells.Item(24, 4).Value = "HighProbConvergence"
$wsHPC.Cells.Item(24, 5).Value = "UP"
$wsHPC.Cells.Item(24, 6).Value = 0.96
$wsHPC.Cells.Item(24, 8).Value = "OPEN"
$wsHPC.Cells.Item(24, 9).Value = 0
$wsHPC.Cells.Item(24, 10).Value = 0
$wsHPC.Cells.Item(24, 11).Value = 100.4130057263667
$wsHPC.Cells.Item(24, 12).Value = 0
$wsHPC.Cells.Item(24, 13).Value = 0
$wsHPC.Cells.Item(24, 14).Value = 0.95
$wsHPC.Cells.Item(24, 15).Value = "Mean reversion UP: price 1.75% below mean (z=-3.00)"
$wsHPC.Cells.Item(24, 17).Value = 0

# ---------------------------------------------------------------------------
# 6) MarketMaking sheet
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# 6a) Close trade #132 (row 53, DOWN)
$wsMM.Cells.Item(53, 7).Value = 0.03                # Exit Price
$wsMM.Cells.Item(53, 8).Value = "CLOSED"            # Status
$wsMM.Cells.Item(53, 9).Value = -25                 # P&L %
$wsMM.Cells.Item(53, 10).Value = -0.01              # P&L $
$wsMM.Cells.Item(53, 11).Value = 99.14              # Capital After
$wsMM.Cells.Item(53, 16).Value = "early_exit"       # Exit Reason
$wsMM.Cells.Item(53, 17).Value = 0.15               # Duration (min)

# 6b) New row 64 - trade #163 (UP, OPEN)
$wsMM.Cells.Item(64, 1).Value = 163
Set-TextValue $wsMM.Cells.Item(64, 2) "2026-02-18"
Set-TextValue $wsMM.Cells.Item(64, 3) "00:37:35"
$wsMM.Cells.Item(64, 4).Value = "MarketMaking"
$wsMM.Cells.Item(64, 5).Value = "UP"
$wsMM.Cells.Item(64, 6).Value = 0.9399999999999999
$wsMM.Cells.Item(64, 8).Value = "OPEN"
$wsMM.Cells.Item(64, 9).Value = 0
$wsMM.Cells.Item(64, 10).Value = 0
$wsMM.Cells.Item(64, 11).Value = 99.14858346467945
$wsMM.Cells.Item(64, 12).Value = 0
$wsMM.Cells.Item(64, 13).Value = 0
$wsMM.Cells.Item(64, 14).Value = 0.6
$wsMM.Cells.Item(64, 15).Value = "Normal spread capture: 198 bps"
$wsMM.Cells.Item(64, 17).Value = 0
